$d = $word.ActiveDocument

# 1. Update the letter date (unique occurrence in the document).
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the single-line mailing address "969 Story Road, San Jose CA 95122"
#    into two separate paragraphs: "969 Story Road" and "San Jose, CA 95122".
#    NOTE: this exact string also appears later inside the "PROPERTY ADDRESS"
#    table, which must stay untouched, so restrict the search to the range
#    that precedes the document's first table (the letterhead block).
$tableStart = $d.Tables(1).Range.Start
$addrScope = $d.Range(0, $tableStart)
$addrScope.Find.Execute("969 Story Road, San Jose CA 95122", $true, $false, $false, $false, $false, $true, 1, $false, "969 Story RoadSan Jose, CA 95122", 2) | Out-Null

$addrScope2 = $d.Range(0, $tableStart)
$addrScope2.Find.Execute("969 Story RoadSan Jose, CA 95122") | Out-Null
$splitPos = $addrScope2.Start + ("969 Story Road").Length
$insertRange = $d.Range($splitPos, $splitPos)
$insertRange.InsertParagraphAfter()

# 3. Remove the empty "No Spacing" paragraph that immediately follows the
#    "Board of Directors" signature line.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs($i)
  if ($p.Range.Text -match "Board of Directors") {
    $next = $d.Paragraphs($i + 1)
    if ($next.Style.NameLocal -eq "No Spacing" -and $next.Range.Text.Trim() -eq "") {
      $next.Range.Delete()
    }
    break
  }
}
